$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 25,9
$data[0,0] = "model_6_9_6"
$data[0,1] = 0.2304671916459294
$data[0,2] = 0.1768808165800588
$data[0,3] = -0.04420886564384485
$data[0,4] = 0.08688148562427289
$data[0,5] = 0.8516460657119751
$data[0,6] = 1.420858263969421
$data[0,7] = 1.736180067062378
$data[0,8] = 1.569245219230652
$data[1,0] = "model_6_9_5"
$data[1,1] = 0.2312253711294555
$data[1,2] = 0.1812789445634206
$data[1,3] = -0.03507737528146171
$data[1,4] = 0.09337827305343283
$data[1,5] = 0.8508070111274719
$data[1,6] = 1.413266181945801
$data[1,7] = 1.720997333526611
$data[1,8] = 1.558079957962036
$data[2,0] = "model_6_9_4"
$data[2,1] = 0.2375054988532181
$data[2,2] = 0.2064208633940718
$data[2,3] = -0.01916077088540757
$data[2,4] = 0.113993212236509
$data[2,5] = 0.8438567519187927
$data[2,6] = 1.369866609573364
$data[2,7] = 1.694533109664917
$data[2,8] = 1.52265202999115
$data[3,0] = "model_6_9_7"
$data[3,1] = 0.2386122376102731
$data[3,2] = 0.1903529386520749
$data[3,3] = -0.03040041255157755
$data[3,4] = 0.1003317099604351
$data[3,5] = 0.8426318764686584
$data[3,6] = 1.397602796554565
$data[3,7] = 1.713220953941345
$data[3,8] = 1.546130061149597
$data[4,0] = "model_6_9_24"
$data[4,1] = 0.2447185745659718
$data[4,2] = 0.1965481154384218
$data[4,3] = -0.01814805844272427
$data[4,4] = 0.1092051759614291
$data[4,5] = 0.8358739614486694
$data[4,6] = 1.386908769607544
$data[4,7] = 1.692849397659302
$data[4,8] = 1.530880451202393
$data[5,0] = "model_6_9_23"
$data[5,1] = 0.2448761140652409
$data[5,2] = 0.1964323456909794
$data[5,3] = -0.01724626961785591
$data[5,4] = 0.1095523574561462
$data[5,5] = 0.8356996178627014
$data[5,6] = 1.387108683586121
$data[5,7] = 1.691349983215332
$data[5,8] = 1.53028392791748
$data[6,0] = "model_6_9_21"
$data[6,1] = 0.244942555645197
$data[6,2] = 0.1959589690371648
$data[6,3] = -0.01628565271234672
$data[6,4] = 0.1097395418563395
$data[6,5] = 0.8356261253356934
$data[6,6] = 1.387925744056702
$data[6,7] = 1.689752817153931
$data[6,8] = 1.529962182044983
$data[7,0] = "model_6_9_22"
$data[7,1] = 0.244957812137559
$data[7,2] = 0.1963262724411134
$data[7,3] = -0.01670416060518476
$data[7,4] = 0.1097445670310526
$data[7,5] = 0.8356092572212219
$data[7,6] = 1.387291669845581
$data[7,7] = 1.690448522567749
$data[7,8] = 1.529953718185425
$data[8,0] = "model_6_9_20"
$data[8,1] = 0.2459851130459003
$data[8,2] = 0.1998225203991028
$data[8,3] = -0.0158157902535101
$data[8,4] = 0.112006434716584
$data[8,5] = 0.8344722986221313
$data[8,6] = 1.381256461143494
$data[8,7] = 1.688971519470215
$data[8,8] = 1.526066541671753
$data[9,0] = "model_6_9_17"
$data[9,1] = 0.2463144344173652
$data[9,2] = 0.2003403041868755
$data[9,3] = -0.0145489225970401
$data[9,4] = 0.1128601625951433
$data[9,5] = 0.8341078758239746
$data[9,6] = 1.380362749099731
$data[9,7] = 1.686865091323853
$data[9,8] = 1.524599194526672
$data[10,0] = "model_6_9_16"
$data[10,1] = 0.2463763002364324
$data[10,2] = 0.2004119249671358
$data[10,3] = -0.01427217060682895
$data[10,4] = 0.1130243338581149
$data[10,5] = 0.8340393900871277
$data[10,6] = 1.380239009857178
$data[10,7] = 1.686404943466187
$data[10,8] = 1.524317145347595
$data[11,0] = "model_6_9_19"
$data[11,1] = 0.2465027709634792
$data[11,2] = 0.2009400648713141
$data[11,3] = -0.01460456569716362
$data[11,4] = 0.1131538639162211
$data[11,5] = 0.8338994383811951
$data[11,6] = 1.379327416419983
$data[11,7] = 1.686957597732544
$data[11,8] = 1.524094581604004
$data[12,0] = "model_6_9_14"
$data[12,1] = 0.246547375131117
$data[12,2] = 0.2007666843863244
$data[12,3] = -0.01329349373889022
$data[12,4] = 0.1136568606644398
$data[12,5] = 0.8338499665260315
$data[12,6] = 1.379626870155334
$data[12,7] = 1.684777736663818
$data[12,8] = 1.523230075836182
$data[13,0] = "model_6_9_18"
$data[13,1] = 0.2465857823437831
$data[13,2] = 0.2008554402768146
$data[13,3] = -0.01389551187538207
$data[13,4] = 0.1134315532356887
$data[13,5] = 0.8338075280189514
$data[13,6] = 1.379473447799683
$data[13,7] = 1.685778737068176
$data[13,8] = 1.523617386817932
$data[14,0] = "model_6_9_8"
$data[14,1] = 0.2469330739962197
$data[14,2] = 0.2062586041403173
$data[14,3] = -0.01374505112369162
$data[14,4] = 0.1163728089010269
$data[14,5] = 0.8334231972694397
$data[14,6] = 1.370146632194519
$data[14,7] = 1.685528516769409
$data[14,8] = 1.51856255531311
$data[15,0] = "model_6_9_15"
$data[15,1] = 0.2474012226905242
$data[15,2] = 0.2030929366147158
$data[15,3] = -0.01246950895297316
$data[15,4] = 0.1152706736812906
$data[15,5] = 0.832905113697052
$data[15,6] = 1.375611186027527
$data[15,7] = 1.683407664299011
$data[15,8] = 1.520456790924072
$data[16,0] = "model_6_9_13"
$data[16,1] = 0.2475305910927309
$data[16,2] = 0.2037685088805898
$data[16,3] = -0.01221762055329956
$data[16,4] = 0.1157447966895644
$data[16,5] = 0.8327618837356567
$data[16,6] = 1.374445080757141
$data[16,7] = 1.682988882064819
$data[16,8] = 1.519641876220703
$data[17,0] = "model_6_9_12"
$data[17,1] = 0.2475482963712199
$data[17,2] = 0.2047267224064254
$data[17,3] = -0.01291860463166561
$data[17,4] = 0.1159350010696888
$data[17,5] = 0.8327422738075256
$data[17,6] = 1.372790932655334
$data[17,7] = 1.684154510498047
$data[17,8] = 1.519314885139465
$data[18,0] = "model_6_9_9"
$data[18,1] = 0.2490782081484925
$data[18,2] = 0.2101481564039916
$data[18,3] = -0.009985159103707231
$data[18,4] = 0.1201534330854377
$data[18,5] = 0.8310491442680359
$data[18,6] = 1.36343252658844
$data[18,7] = 1.679277062416077
$data[18,8] = 1.512065291404724
$data[19,0] = "model_6_9_10"
$data[19,1] = 0.2491347606798847
$data[19,2] = 0.2103452026468442
$data[19,3] = -0.01126016228445526
$data[19,4] = 0.1196777627068077
$data[19,5] = 0.8309865593910217
$data[19,6] = 1.363092541694641
$data[19,7] = 1.681396961212158
$data[19,8] = 1.512882828712463
$data[20,0] = "model_6_9_11"
$data[20,1] = 0.2491794903960163
$data[20,2] = 0.2104448138404922
$data[20,3] = -0.01122819740653269
$data[20,4] = 0.1197459418427569
$data[20,5] = 0.8309370875358582
$data[20,6] = 1.36292040348053
$data[20,7] = 1.681343793869019
$data[20,8] = 1.512765645980835
$data[21,0] = "model_6_9_3"
$data[21,1] = 0.2839872291672964
$data[21,2] = 0.3043525491180613
$data[21,3] = 0.09885820857994332
$data[21,4] = 0.219802800872584
$data[21,5] = 0.7924151420593262
$data[21,6] = 1.200818061828613
$data[21,7] = 1.498306035995483
$data[21,8] = 1.3408123254776
$data[22,0] = "model_6_9_0"
$data[22,1] = 0.2937639237727891
$data[22,2] = 0.429562355835966
$data[22,3] = 0.4423156606101083
$data[22,4] = 0.4427549617416759
$data[22,5] = 0.7815952301025391
$data[22,6] = 0.9846825003623962
$data[22,7] = 0.9272478222846985
$data[22,8] = 0.9576566219329834
$data[23,0] = "model_6_9_1"
$data[23,1] = 0.294081177892726
$data[23,2] = 0.4285824535620977
$data[23,3] = 0.4409662566835698
$data[23,4] = 0.4416207748916733
$data[23,5] = 0.781244158744812
$data[23,6] = 0.9863739013671875
$data[23,7] = 0.929491400718689
$data[23,8] = 0.9596056938171387
$data[24,0] = "model_6_9_2"
$data[24,1] = 0.3096338858669928
$data[24,2] = 0.3791921348033127
$data[24,3] = 0.2456399293095876
$data[24,4] = 0.3264273092451718
$data[24,5] = 0.7640318274497986
$data[24,6] = 1.071630835533142
$data[24,7] = 1.254255652427673
$data[24,8] = 1.157572269439697

$ws.Range("A2:I26").Value = $data
